$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg).
# This is the weekly reshuffle of the "Fruta, Femacal de La Calera - Coco" data.
$rows = @{
    2  = @{ D = 44403; M = 50; N = 20000; O = 20000; P = 20000; S = 1000 }
    3  = @{ D = 44448; M = 30; N = 22000; O = 22000; P = 22000; S = 1100 }
    4  = @{ D = 44406; M = 20; N = 20000; O = 20000; P = 20000; S = 1000 }
    5  = @{ D = 44291; M = 70; N = 25000; O = 25000; P = 25000; S = 1250 }
    6  = @{ D = 44400; M = 45; N = 20000; O = 20000; P = 20000; S = 1000 }
    7  = @{ D = 44307; M = 30; N = 22000; O = 22000; P = 22000; S = 1100 }
    8  = @{ D = 44298; M = 65; N = 22000; O = 22000; P = 22000; S = 1100 }
    9  = @{ D = 44300; M = 45; N = 22000; O = 22000; P = 22000; S = 1100 }
    10 = @{ D = 44301; M = 38; N = 22000; O = 22000; P = 22000; S = 1100 }
    11 = @{ D = 44377; M = 25; N = 20000; O = 20000; P = 20000; S = 1000 }
    12 = @{ D = 44292; M = 30; N = 25000; O = 25000; P = 25000; S = 1250 }
    13 = @{ D = 44413; M = 45; N = 20000; O = 20000; P = 20000; S = 1000 }
    14 = @{ D = 44382; M = 24; N = 20000; O = 20000; P = 20000; S = 1000 }
    15 = @{ D = 44305; M = 20; N = 22000; O = 22000; P = 22000; S = 1100 }
    16 = @{ D = 44445; M = 45; N = 20000; O = 20000; P = 20000; S = 1000 }
    17 = @{ D = 44376; M = 38; N = 20000; O = 20000; P = 20000; S = 1000 }
    18 = @{ D = 44389; M = 20; N = 20000; O = 20000; P = 20000; S = 1000 }
    19 = @{ D = 44294; M = 25; N = 25000; O = 25000; P = 25000; S = 1250 }
    20 = @{ D = 44385; M = 36; N = 20000; O = 20000; P = 20000; S = 1000 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("S$r").Value = $vals.S
}
